# Auto-generated cell updates from the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.812.31"
$ws.Range("E2").Value = "'  +8.62%  "
$ws.Range("D3").Value = "'1.952.98"
$ws.Range("E3").Value = "'  +6.73%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "'  -0.22%  "
$ws.Range("D5").Value = "'342.70"
$ws.Range("E5").Value = "'  +3.45%  "
$ws.Range("D6").Value = "'0.9989"
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("D7").Value = "'0.4797"
$ws.Range("E7").Value = "'  +4.62%  "
$ws.Range("D8").Value = "'0.4148"
$ws.Range("E8").Value = "'  +8.52%  "
$ws.Range("D9").Value = "'48.55"
$ws.Range("E9").Value = "'  +4.55%  "
$ws.Range("D10").Value = "'0.08294"
$ws.Range("E10").Value = "'  +5.12%  "
$ws.Range("D11").Value = "'1.054"
$ws.Range("E11").Value = "'  +8.84%  "
$ws.Range("D12").Value = "'22.92"
$ws.Range("E12").Value = "'  +8.55%  "
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.941.54"
$ws.Range("E13").Value = "'  +6.00%  "
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.199"
$ws.Range("E14").Value = "'  +5.42%  "
$ws.Range("D15").Value = "'7.470"
$ws.Range("E15").Value = "'  +5.30%  "
$ws.Range("D16").Value = "'93.30"
$ws.Range("E16").Value = "'  +4.07%  "
$ws.Range("D17").Value = "'0.9994"
$ws.Range("E17").Value = "'  -0.22%  "
$ws.Range("E18").Value = "'  +4.40%  "
$ws.Range("D19").Value = "'0.06696"
$ws.Range("E19").Value = "'  +1.50%  "
$ws.Range("D20").Value = "'18.11"
$ws.Range("E20").Value = "'  +5.35%  "
$ws.Range("D21").Value = "'0.9972"
$ws.Range("E21").Value = "'  -0.39%  "
$ws.Range("D22").Value = "'29.770.99"
$ws.Range("E22").Value = "'  +8.51%  "
$ws.Range("D23").Value = "'5.641"
$ws.Range("E23").Value = "'  +5.89%  "
$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "'  +4.68%  "
$ws.Range("D25").Value = "'2.274"
$ws.Range("E25").Value = "'  -0.27%  "
$ws.Range("D26").Value = "'2.161.16"
$ws.Range("E26").Value = "'  +5.53%  "
$ws.Range("D27").Value = "'162.09"
$ws.Range("E27").Value = "'  +3.94%  "
$ws.Range("D28").Value = "'20.24"
$ws.Range("E28").Value = "'  +4.27%  "
$ws.Range("D29").Value = "'2.211"
$ws.Range("E29").Value = "'  +6.69%  "
$ws.Range("D30").Value = "'5.687"
$ws.Range("E30").Value = "'  +7.34%  "
$ws.Range("D31").Value = "'122.88"
$ws.Range("E31").Value = "'  +3.82%  "
$ws.Range("D32").Value = "'1.032"
$ws.Range("E32").Value = "'  +9.88%  "
$ws.Range("D33").Value = "'0.09664"
$ws.Range("E33").Value = "'  +3.87%  "
$ws.Range("E34").Value = "'  +11.57%  "
$ws.Range("D35").Value = "'3.679"
$ws.Range("E35").Value = "'  +2.93%  "
$ws.Range("D36").Value = "'5.519"
$ws.Range("D37").Value = "'0.06274"
$ws.Range("E37").Value = "'  +5.67%  "
$ws.Range("D38").Value = "'0.02322"
$ws.Range("E38").Value = "'  +6.62%  "
$ws.Range("D39").Value = "'8.746"
$ws.Range("E39").Value = "'  +7.51%  "
$ws.Range("D40").Value = "'1.210"
$ws.Range("E40").Value = "'  +5.92%  "
$ws.Range("D41").Value = "'0.6134"
$ws.Range("E41").Value = "'  +6.19%  "
$ws.Range("E42").Value = "'  +7.94%  "
$ws.Range("D43").Value = "'0.1915"
$ws.Range("E43").Value = "'  +4.70%  "
$ws.Range("D44").Value = "'0.9985"
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("D45").Value = "'1.297"
$ws.Range("E45").Value = "'  +2.71%  "
$ws.Range("B46").Value = "'Decentraland"
$ws.Range("C46").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5758"
$ws.Range("E46").Value = "'  +5.98%  "
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.59"
$ws.Range("E47").Value = "'  +5.03%  "
$ws.Range("D48").Value = "'2.344"
$ws.Range("E48").Value = "'  +28.00%  "
$ws.Range("E49").Value = "'  +7.17%  "
$ws.Range("D50").Value = "'0.07316"
$ws.Range("E50").Value = "'  +11.12%  "
$ws.Range("D51").Value = "'114.39"
$ws.Range("E51").Value = "'  +3.96%  "
